$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, pushing existing rows 151-175 down to 152-176.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new weekly record.
$ws.Range("A151").Value = 5
$ws.Range("B151").Value = "Macroferia Regional de Talca"
$ws.Range("C151").Value = "Maule"
$ws.Range("D151").Value = 44522
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100108
$ws.Range("H151").Value = "Tropicales y subtropicales"
$ws.Range("I151").Value = 100108005
$ws.Range("J151").Value = "Piña"
$ws.Range("K151").Value = "Caramelo"
$ws.Range("L151").Value = "Segunda"
$ws.Range("M151").Value = 180
$ws.Range("N151").Value = 19000
$ws.Range("O151").Value = 19000
$ws.Range("P151").Value = 19000
$ws.Range("Q151").Value = "$/caja 14 unidades"
$ws.Range("R151").Value = "Ecuador"
$ws.Range("S151").Value = 1357
$ws.Range("T151").Value = 14
